$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 366, shifting the existing rows
# 366-481 down to 367-482 (dimension grows from A1:R481 to A1:R482).
$ws.Rows("366:366").Insert()

# Populate the newly inserted row 366 with its data.
$ws.Range("A366").Value = 9
$ws.Range("B366").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C366").Value = "Metropolitana"
$ws.Range("D366").Value = 44876
$ws.Range("E366").Value = 13
$ws.Range("F366").Value = 100112039
$ws.Range("G366").Value = "Ciboulette"
$ws.Range("H366").Value = "Sin especificar"
$ws.Range("I366").Value = "Primera"
$ws.Range("J366").Value = 620
$ws.Range("K366").Value = 1000
$ws.Range("L366").Value = 1200
$ws.Range("M366").Value = 1097
$ws.Range("N366").Value = "$/docena de atados"
$ws.Range("O366").Value = "Región Metropolitana"
$ws.Range("P366").Value = 366
$ws.Range("Q366").Value = 3
$ws.Range("R366").Value = "Hortaliza"
